$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.735.31"
$ws.Range("E2").Value = "'  +2.92%  "
$ws.Range("D3").Value = "'3.321.52"
$ws.Range("E3").Value = "'  +0.56%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("D5").Value = "'587.53"
$ws.Range("E5").Value = "'  +5.41%  "
$ws.Range("D6").Value = "'182.41"
$ws.Range("E6").Value = "'  -0.92%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "'  +3.58%  "
$ws.Range("D9").Value = "'3.313.04"
$ws.Range("E9").Value = "'  +0.64%  "
$ws.Range("E10").Value = "'  +2.73%  "
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "'  +0.93%  "
$ws.Range("D12").Value = "'46.29"
$ws.Range("E12").Value = "'  +1.81%  "
$ws.Range("E13").Value = "'  +5.37%  "
$ws.Range("D14").Value = "'637.81"
$ws.Range("E14").Value = "'  +10.67%  "
$ws.Range("D15").Value = "'3.857.78"
$ws.Range("E15").Value = "'  +0.69%  "
$ws.Range("D16").Value = "'8.44"
$ws.Range("E16").Value = "'  +0.95%  "
$ws.Range("D17").Value = "'67.883.60"
$ws.Range("E17").Value = "'  +3.29%  "
$ws.Range("E18").Value = "'  +1.73%  "
$ws.Range("D19").Value = "'3.328.07"
$ws.Range("E19").Value = "'  +0.64%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("E20").Value = "'  +1.10%  "
$ws.Range("D21").Value = "'10.89"
$ws.Range("E21").Value = "'  +1.07%  "
$ws.Range("D22").Value = "'0.900"
$ws.Range("E22").Value = "'  +1.77%  "
$ws.Range("D23").Value = "'17.65"
$ws.Range("E23").Value = "'  -0.24%  "
$ws.Range("D24").Value = "'5.00"
$ws.Range("E24").Value = "'  +0.91%  "
$ws.Range("D25").Value = "'97.28"
$ws.Range("E25").Value = "'  -0.05%  "
$ws.Range("D26").Value = "'4.00"
$ws.Range("E26").Value = "'  +2.25%  "
$ws.Range("E27").Value = "'  +4.11%  "
$ws.Range("D28").Value = "'9.57"
$ws.Range("E28").Value = "'  +3.53%  "
$ws.Range("D29").Value = "'32.81"
$ws.Range("E29").Value = "'  +8.16%  "
$ws.Range("D30").Value = "'8.56"
$ws.Range("E30").Value = "'  +2.62%  "
$ws.Range("E31").Value = "'  +1.20%  "
$ws.Range("D32").Value = "'593.72"
$ws.Range("E32").Value = "'  +6.08%  "
$ws.Range("D33").Value = "'3.934.12"
$ws.Range("E33").Value = "'  +5.71%  "
$ws.Range("E34").Value = "'  +2.04%  "
$ws.Range("D35").Value = "'3.52"
$ws.Range("E35").Value = "'  -4.06%  "
$ws.Range("E36").Value = "'  +1.93%  "
$ws.Range("E37").Value = "'  -0.18%  "
$ws.Range("E38").Value = "'  +0.36%  "
$ws.Range("E39").Value = "'  +4.48%  "
$ws.Range("D40").Value = "'0.129"
$ws.Range("E40").Value = "'  +1.98%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "'  +4.70%  "
$ws.Range("D42").Value = "'32.61"
$ws.Range("E42").Value = "'  -0.99%  "
$ws.Range("B43").Value = "'PEPE"
$ws.Range("C43").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "'0.0₃0683"
$ws.Range("E43").Value = "'  +1.35%  "
$ws.Range("B44").Value = "'ApeXProtocol"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.38"
$ws.Range("E44").Value = "'  +1.44%  "
$ws.Range("E45").Value = "'  +2.40%  "
$ws.Range("E46").Value = "'  +2.32%  "
$ws.Range("E47").Value = "'  +2.24%  "
$ws.Range("E48").Value = "'  +0.68%  "
$ws.Range("E49").Value = "'  +2.46%  "
$ws.Range("E50").Value = "'  +9.96%  "
$ws.Range("D51").Value = "'130.82"
$ws.Range("E51").Value = "'  +5.11%  "
